$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("D2", "E2", "D3", "E3", "D4", "E4", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "E20", "E21", "E22", "D23", "E23", "E24", "E25", "D26", "E26", "E27", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "D49", "E49", "D50", "E50")
foreach ($cellAddr in $targetCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "261.62"
$ws.Range("E2").Value = "1.70%"
$ws.Range("D3").Value = "27.23"
$ws.Range("E3").Value = "0.50%"
$ws.Range("D4").Value = "4.719"
$ws.Range("E4").Value = "2.83%"
$ws.Range("E5").Value = "2.98%"
$ws.Range("D6").Value = "6.660"
$ws.Range("E6").Value = "0.39%"
$ws.Range("D7").Value = "0.8632"
$ws.Range("E7").Value = "1.36%"
$ws.Range("D8").Value = "0.9226"
$ws.Range("E8").Value = "-2.12%"
$ws.Range("D9").Value = "0.1410"
$ws.Range("E9").Value = "1.18%"
$ws.Range("D10").Value = "0.05101"
$ws.Range("E10").Value = "2.65%"
$ws.Range("D11").Value = "0.07107"
$ws.Range("E11").Value = "0.45%"
$ws.Range("D12").Value = "0.03048"
$ws.Range("E12").Value = "-1.15%"
$ws.Range("D13").Value = "0.09096"
$ws.Range("E13").Value = "-0.42%"
$ws.Range("D14").Value = "0.001533"
$ws.Range("E14").Value = "0.49%"
$ws.Range("D15").Value = "0.0006100"
$ws.Range("E15").Value = "0.93%"
$ws.Range("D16").Value = "0.006046"
$ws.Range("E16").Value = "-0.89%"
$ws.Range("D17").Value = "3.446"
$ws.Range("E17").Value = "-1.40%"
$ws.Range("E18").Value = "-0.44%"
$ws.Range("E20").Value = "2.40%"
$ws.Range("E21").Value = "2.15%"
$ws.Range("E22").Value = "3.77%"
$ws.Range("D23").Value = "0.04241"
$ws.Range("E23").Value = "-0.45%"
$ws.Range("E24").Value = "-0.30%"
$ws.Range("E25").Value = "-8.72%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.02%"
$ws.Range("E27").Value = "-18.95%"
$ws.Range("D40").Value = "0.03881"
$ws.Range("E40").Value = "1.58%"
$ws.Range("D41").Value = "0.1113"
$ws.Range("E41").Value = "1.10%"
$ws.Range("E42").Value = "-34.13%"
$ws.Range("D43").Value = "0.01501"
$ws.Range("E43").Value = "5.98%"
$ws.Range("D44").Value = "0.002182"
$ws.Range("E44").Value = "-9.49%"
$ws.Range("D45").Value = "0.00005322"
$ws.Range("E45").Value = "-1.04%"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("E47").Value = "6.92%"
$ws.Range("D48").Value = "0.1353"
$ws.Range("E48").Value = "-46.24%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "-0.05%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "-0.05%"
